# "Generate Report for Archive"
#
# The localization status report is regenerated: the status text for the
# one row that was still mid-flight moves from "Ready for handoff" to
# "In Translation" on all three sheets (Overview!E2/F2, zh-cn!C2,
# de-de!C2 all shared the same string), and the status column on each
# sheet is re-sized to fit the new (shorter) text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Sheets.Item("Overview")
$zhcn     = $wb.Sheets.Item("zh-cn")
$dede     = $wb.Sheets.Item("de-de")

# Update the status text everywhere it appears.
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# Re-fit the status columns now that the text is shorter.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
